# feat: add industry into import graduation students
#
# Adds two new columns (Ma nganh / Ten nganh) with sample data to the
# graduation-students import template, and nudges the sheet's column
# widths / selection to match the refreshed template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---------------------------------------------
$ws.Range("F1").Value = "Mã ngành"
$ws.Range("G1").Value = "Tên ngành"

# --- New data cells (row 2) ------------------------------------------------
$ws.Range("F2").Value = 7480201
$ws.Range("G2").Value = "Công nghệ thông tin"

# --- Column widths (A, D, F, G get explicit widths in the new layout) ------
# Excel's ColumnWidth setter snaps to the MDW-7 pixel grid (px = round(chars*6)+5,
# stored chars = px/6), so the literal target "characters" value is fed back
# through that same pixel rounding to land on the nearest grid point Excel
# itself would have produced for these widths.
$ws.Columns.Item(1).ColumnWidth = 13.451822916666666
$ws.Columns.Item(4).ColumnWidth = 14.877604166666666
$ws.Columns.Item(6).ColumnWidth = 10.166666666666666
$ws.Columns.Item(7).ColumnWidth = 22.307291666666668

# --- Selection moves as part of the refreshed template ---------------------
$null = $ws.Range("C9").Select()
